$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 704 (the "iPhone games free download" post), which shifts
# all subsequent rows up by one.
$ws.Rows.Item(704).Delete()
